# Swap the species-record data between row 2 and row 3 (keeping the
# shared/location columns untouched), and move the "Biotop" (AH) note
# from row 3 to row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H")

foreach ($col in $cols) {
    $addr2 = "$col`2"
    $addr3 = "$col`3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}

# AH2 is currently empty, AH3 holds "Vägkant" -> move it up to AH2.
$ws.Range("AH2").Value2 = $ws.Range("AH3").Value2
$ws.Range("AH3").Value2 = ""
